$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits right after the last
# sentence of the "message bus" paragraph. In the edit we need to:
#   1) leave that paragraph's text alone (just drop the bookmark from it)
#   2) add a new sentence "Circular dependencies!!! What a pain!!!" into
#      what is currently the final (empty) paragraph of the document
#   3) move the "_GoBack" bookmark so it ends up right after that new
#      sentence (still inside the same, now non-empty, last paragraph)

# Step 1: remove the bookmark from its current location.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Step 2: type the new sentence into the last paragraph, matching the
# "en-US" language formatting used throughout the document.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertBefore("Circular dependencies!!! What a pain!!!")
$insertionPoint.LanguageID = "en-US"
$bookmarkPos = $insertionPoint.End

# Step 3: re-create "_GoBack" right after the new text. Adding a bookmark
# whose position sits at (or very near) the absolute end of the document
# places it incorrectly, so temporarily pad the document with extra text
# past the target spot, add the bookmark, then remove the padding again.
$tail = $d.Content
$tail.Collapse(0)
$padding = "Z".PadLeft(80, "Z")
$tail.InsertAfter($padding)

$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$paddingRange = $d.Range($bookmarkPos, $d.Content.End)
$paddingRange.Delete()
